$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2892
$ws.Range("L3").Value = 2929
$ws.Range("L4").Value = 777
$ws.Range("L6").Value = 2630
$ws.Range("L7").Value = 9393

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 74
$ws.Range("L7").Value = 315
$ws.Range("L8").Value = 594
$ws.Range("L10").Value = 60
$ws.Range("L14").Value = 45
$ws.Range("L15").Value = 69
$ws.Range("L19").Value = 267
$ws.Range("L27").Value = 90
$ws.Range("L29").Value = 503
$ws.Range("L31").Value = 90
$ws.Range("L33").Value = 429
$ws.Range("L34").Value = 56
$ws.Range("L36").Value = 130
$ws.Range("L37").Value = 344
$ws.Range("L42").Value = 305
$ws.Range("L44").Value = 72
$ws.Range("L45").Value = 17
$ws.Range("L47").Value = 74
$ws.Range("L48").Value = 125
$ws.Range("L51").Value = 114
$ws.Range("L52").Value = 189
$ws.Range("L63").Value = 30
$ws.Range("L64").Value = 58
$ws.Range("L65").Value = 168
$ws.Range("L67").Value = 348
$ws.Range("L71").Value = 28
$ws.Range("L72").Value = 44
$ws.Range("L79").Value = 251
$ws.Range("L83").Value = 221
$ws.Range("L85").Value = 478
$ws.Range("L86").Value = 70
$ws.Range("L87").Value = 30
$ws.Range("L88").Value = 118
$ws.Range("L89").Value = 121
$ws.Range("L90").Value = 92
$ws.Range("L91").Value = 133
$ws.Range("L94").Value = 114
$ws.Range("L95").Value = 125
$ws.Range("L96").Value = 94
$ws.Range("L99").Value = 156
$ws.Range("L101").Value = 9393

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 19
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 96
$ws.Range("L3").Value = 99
$ws.Range("L4").Value = 25
$ws.Range("L7").Value = 315

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 121

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 194
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 478

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 64
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 200
$ws.Range("L7").Value = 594

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 221

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 121
$ws.Range("L6").Value = 148
$ws.Range("L7").Value = 429

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 102
$ws.Range("L3").Value = 100
$ws.Range("L6").Value = 109
$ws.Range("L7").Value = 344

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 54
$ws.Range("L7").Value = 168

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 39
$ws.Range("L3").Value = 63
$ws.Range("L7").Value = 156

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 128
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 348

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 134
$ws.Range("L7").Value = 503

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 93
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 305

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 16
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 89
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 31
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 39
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 33
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 30

Write-Host "Applied 137 cell updates for 2025-06-14 data"